$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 77, shifting existing rows 77:105 down to 78:106
# (mirrors the weekly data refresh: a new Jengibre price observation is
# prepended to the series and all subsequent rows move down by one)
$ws.Rows(77).Insert()

# Populate the newly inserted row 77 with the new weekly observation
$ws.Range("A77").Value = 9
$ws.Range("B77").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C77").Value = "Metropolitana"
$ws.Range("D77").Value = 44795
$ws.Range("E77").Value = 13
$ws.Range("F77").Value = 100114007
$ws.Range("G77").Value = "Jengibre"
$ws.Range("H77").Value = "Sin especificar"
$ws.Range("I77").Value = "Primera"
$ws.Range("J77").Value = 700
$ws.Range("K77").Value = 15000
$ws.Range("L77").Value = 16000
$ws.Range("M77").Value = 15429
$ws.Range("N77").Value = "`$/caja 13 kilos"
$ws.Range("O77").Value = "Perú"
$ws.Range("P77").Value = 1187
$ws.Range("Q77").Value = 13
$ws.Range("R77").Value = "Hortaliza"
